$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 3.75
$ws.Range("I2").Value = 2.2
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 3.7
$ws.Range("V2").Value = 1.83
$ws.Range("W2").Value = 1.36
$ws.Range("X2").Value = 13.5
$ws.Range("AL2").Value = 60
$ws.Range("O3").Value = 1.17
$ws.Range("Q3").Value = 1.17
$ws.Range("S3").Value = 1.17
$ws.Range("F4").Value = 4.8
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 1.69
$ws.Range("I4").Value = 1.83
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 4.4
$ws.Range("L4").Value = 1.37
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 3.4
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 1.89
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.37
$ws.Range("S4").Value = 2.84
$ws.Range("T4").Value = 1.84
$ws.Range("U4").Value = 1.98
$ws.Range("V4").Value = 2.2
$ws.Range("W4").Value = 1.18
$ws.Range("Y4").Value = 10.5
$ws.Range("AC4").Value = 11.5
$ws.Range("AO4").Value = 11.5
$ws.Range("T5").Value = 1.83
$ws.Range("F6").Value = 2.76
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 1.96
$ws.Range("I6").Value = 2.64
$ws.Range("J6").Value = 2.7
$ws.Range("K6").Value = 8
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 3.8
$ws.Range("P6").Value = 2.06
$ws.Range("Q6").Value = 1.55
$ws.Range("R6").Value = 1.48
$ws.Range("S6").Value = 2.38
$ws.Range("T6").Value = 1.66
$ws.Range("U6").Value = 2.32
$ws.Range("V6").Value = 1.61
$ws.Range("W6").Value = 1.33
$ws.Range("Y6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("F7").Value = 6.4
$ws.Range("G7").Value = 10.5
$ws.Range("H7").Value = 1.41
$ws.Range("I7").Value = 1.53
$ws.Range("J7").Value = 4.5
$ws.Range("K7").Value = 6.4
$ws.Range("L7").Value = 1.01
$ws.Range("N7").Value = 3.8
$ws.Range("O7").Value = 1.24
$ws.Range("P7").Value = 2
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.41
$ws.Range("S7").Value = 2.8
$ws.Range("T7").Value = 1.98
$ws.Range("U7").Value = 1.79
$ws.Range("V7").Value = 2.88
$ws.Range("W7").Value = 1.11
$ws.Range("AF7").Value = 85
$ws.Range("J8").Value = 4
$ws.Range("AC8").Value = 11
$ws.Range("I9").Value = 3.75
$ws.Range("J9").Value = 3.4
$ws.Range("N9").Value = 5.2
$ws.Range("V9").Value = 1.37
$ws.Range("I10").Value = 1.73
$ws.Range("V10").Value = 2.36
$ws.Range("F11").Value = 12.5
$ws.Range("G11").Value = 16
$ws.Range("H11").Value = 1.24
$ws.Range("I11").Value = 1.29
$ws.Range("K11").Value = 7.8
$ws.Range("N11").Value = 8.2
$ws.Range("O11").Value = 1.11
$ws.Range("P11").Value = 3.4
$ws.Range("Q11").Value = 1.33
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.82
$ws.Range("T11").Value = 1.73
$ws.Range("U11").Value = 2.1
$ws.Range("V11").Value = 4.4
$ws.Range("W11").Value = 1.06
$ws.Range("X11").Value = 48
$ws.Range("Y11").Value = 16
$ws.Range("Z11").Value = 12
$ws.Range("AD11").Value = 12.5
$ws.Range("AF11").Value = 180
$ws.Range("AJ11").Value = 510
$ws.Range("AK11").Value = 210
$ws.Range("AL11").Value = 150
$ws.Range("F12").Value = 4.6
$ws.Range("G12").Value = 5.7
$ws.Range("H12").Value = 1.68
$ws.Range("I12").Value = 1.81
$ws.Range("J12").Value = 4
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 4.8
$ws.Range("O12").Value = 1.19
$ws.Range("P12").Value = 2.3
$ws.Range("Q12").Value = 1.62
$ws.Range("R12").Value = 1.52
$ws.Range("S12").Value = 2.46
$ws.Range("T12").Value = 1.65
$ws.Range("U12").Value = 2.22
$ws.Range("V12").Value = 2.22
$ws.Range("W12").Value = 1.22
$ws.Range("X12").Value = 27
$ws.Range("Y12").Value = 13.5
$ws.Range("Z12").Value = 15
$ws.Range("AA12").Value = 22
$ws.Range("AB12").Value = 27
$ws.Range("AC12").Value = 12.5
$ws.Range("AD12").Value = 12.5
$ws.Range("AE12").Value = 21
$ws.Range("AF12").Value = 50
$ws.Range("AG12").Value = 24
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 36
$ws.Range("AK12").Value = 70
$ws.Range("AL12").Value = 70
$ws.Range("AM12").Value = 95
$ws.Range("AN12").Value = 65
$ws.Range("AO12").Value = 9.8
$ws.Range("G13").Value = 1.34
$ws.Range("I13").Value = 14
$ws.Range("J13").Value = 5.4
$ws.Range("P13").Value = 2.6
$ws.Range("R13").Value = 1.64
$ws.Range("V13").Value = 1.07
$ws.Range("W13").Value = 3.9
$ws.Range("AB13").Value = 1000
$ws.Range("AC13").Value = 1000
$ws.Range("AF13").Value = 1000
$ws.Range("AG13").Value = 990
$ws.Range("AJ13").Value = 1000
$ws.Range("AK13").Value = 1000
$ws.Range("AN13").Value = 4.8
